# Adds 3 new columns (CH4, N2O, strøm mean-reduction CO2-eq) to the KVIK emission table.
# Columns E:F:G are inserted before the existing "CH4_dyr_stald" column (old E), shifting
# the existing E..J columns to H..M automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank columns at E:G (existing columns E..J shift right to H..M).
$ws.Range("E1:G1").EntireColumn.Insert()

# New column headers.
$ws.Range("E1").Value = "mean_reduktion_CH4_CO2_eq_tot_m3"
$ws.Range("F1").Value = "mean_reduktion_N2O_CO2_eq_tot_m3"
$ws.Range("G1").Value = "mean_reduktion_strøm_CO2_eq_tot_m3"

# New column values for data rows 2-26.
$efg = New-Object 'object[,]' 25,3
$efg[0,0] = 73.7978673334683; $efg[0,1] = 9.14811135590323; $efg[0,2] = 0
$efg[1,0] = 55.7493691804035; $efg[1,1] = 10.205502623922; $efg[1,2] = 0
$efg[2,0] = 52.9206439812139; $efg[2,1] = 0; $efg[2,2] = 0
$efg[3,0] = 43.1537099295105; $efg[3,1] = 0; $efg[3,2] = 0
$efg[4,0] = 97.7007238992046; $efg[4,1] = 1.48037876756502; $efg[4,2] = 0
$efg[5,0] = 68.4377502242438; $efg[5,1] = 0.441890816153944; $efg[5,2] = 0
$efg[6,0] = 0; $efg[6,1] = 0; $efg[6,2] = 0
$efg[7,0] = 0; $efg[7,1] = 0; $efg[7,2] = 0
$efg[8,0] = 3.17867837967558; $efg[8,1] = 1.76235567567264; $efg[8,2] = -1.48046723981568
$efg[9,0] = 83.6447495475407; $efg[9,1] = 10.9104670315759; $efg[9,2] = -1.48046723981568
$efg[10,0] = 57.8707506590108; $efg[10,1] = 0.184208751303223; $efg[10,2] = 0
$efg[11,0] = 47.1993702354021; $efg[11,1] = 0.225646871976753; $efg[11,2] = 0
$efg[12,0] = 17.0032442509388; $efg[12,1] = 0; $efg[12,2] = 0
$efg[13,0] = 33.1202172700263; $efg[13,1] = 0; $efg[13,2] = 0
$efg[14,0] = 26.971068705944; $efg[14,1] = 0; $efg[14,2] = 0
$efg[15,0] = 9.78245116315753; $efg[15,1] = 0; $efg[15,2] = 0
$efg[16,0] = 19.473049581511; $efg[16,1] = -0.644734075196555; $efg[16,2] = 0
$efg[17,0] = 101.343909794764; $efg[17,1] = 9.13842915566874; $efg[17,2] = 0
$efg[18,0] = 83.6605874196375; $efg[18,1] = 9.56076854872548; $efg[18,2] = 0
$efg[19,0] = 75.5409009168016; $efg[19,1] = 0.150851788419652; $efg[19,2] = 0
$efg[20,0] = 68.3815138706278; $efg[20,1] = -0.644734075196555; $efg[20,2] = 0
$efg[21,0] = 80.8596286846354; $efg[21,1] = 0.184208751303223; $efg[21,2] = 0
$efg[22,0] = 72.9750515095016; $efg[22,1] = -0.422234743611106; $efg[22,2] = 0
$efg[23,0] = 50.4028055866973; $efg[23,1] = 0; $efg[23,2] = 0
$efg[24,0] = 50.0073633151328; $efg[24,1] = -0.644734075196555; $efg[24,2] = 0
$ws.Range("E2:G26").Value = $efg

# Two existing values in column C (sum_reduktion_totCO2_eq_tot) were recalculated.
$ws.Range("C10").Value = 72.1844592942545
$ws.Range("C11").Value = 1152.12848461255
